{"js": "const replacements = [\n  [\"68\u00f74=17, 0\", \"88\u00f79=9, 7\"],\n  [\"85\u00f77=12, 1\", \"15\u00f75=3, 0\"],\n  [\"50\u00f72=25, 0\", \"65\u00f74=16, 1\"],\n  [\"73\u00f72=36, 1\", \"79\u00f75=15, 4\"],\n  [\"92\u00f72=46, 0\", \"59\u00f72=29, 1\"],\n  [\"67\u00f72=33, 1\", \"41\u00f77=5, 6\"],\n  [\"33\u00f76=5, 3\", \"44\u00f74=11, 0\"],\n  [\"37\u00f72=18, 1\", \"13\u00f72=6, 1\"],\n  [\"57\u00f76=9, 3\", \"54\u00f78=6, 6\"],\n  [\"34\u00f72=17, 0\", \"94\u00f77=13, 3\"],\n  [\"93\u00f76=15, 3\", \"77\u00f73=25, 2\"],\n  [\"35\u00f76=5, 5\", \"39\u00f78=4, 7\"],\n  [\"63\u00f74=15, 3\", \"18\u00f73=6, 0\"],\n  [\"77\u00f74=19, 1\", \"41\u00f74=10, 1\"],\n  [\"81\u00f77=11, 4\", \"99\u00f78=12, 3\"],\n  [\"80\u00f73=26, 2\", \"90\u00f78=11, 2\"],\n  [\"66\u00f72=33, 0\", \"30\u00f76=5, 0\"],\n  [\"92\u00f77=13, 1\", \"19\u00f79=2, 1\"],\n  [\"30\u00f74=7, 2\", \"71\u00f77=10, 1\"],\n  [\"90\u00f78=11, 2\", \"61\u00f79=6, 7\"],\n  [\"64\u00f76=10, 4\", \"43\u00f75=8, 3\"],\n  [\"43\u00f78=5, 3\", \"66\u00f79=7, 3\"],\n  [\"67\u00f78=8, 3\", \"79\u00f73=26, 1\"],\n  [\"66\u00f76=11, 0\", \"46\u00f74=11, 2\"],\n  [\"57\u00f73=19, 0\", \"90\u00f74=22, 2\"],\n];\n\n// Phase 1: locate every target range FIRST, while all the \"old\" strings\n// are still unique in the document. This avoids ambiguity that would\n// occur if a later \"new\" value happened to match an earlier \"old\" value\n// that was already replaced (e.g. \"90\\u00f78=11, 2\" appears both as an\n// original value and as a replacement value elsewhere in this document).\nconst found = [];\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false\n  });\n  results.load(\"items\");\n  found.push({ results, newText, oldText });\n}\nawait context.sync();\n\n// Phase 2: apply the replacements now that every range has been captured.\nfor (const { results, newText, oldText } of found) {\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet-style division table: each content row (1, 5, 9, 13, 17)\n# holds 5 practice problems (one per column). Re-point every cell at its\n# expected (row, col) position directly so the edit is unambiguous even\n# though some of the new values coincide with old values used elsewhere\n# in the table (e.g. \"90\u00f78=11, 2\").\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; Old = \"68\u00f74=17, 0\"; New = \"88\u00f79=9, 7\" },\n    @{ Row = 1; Col = 2; Old = \"85\u00f77=12, 1\"; New = \"15\u00f75=3, 0\" },\n    @{ Row = 1; Col = 3; Old = \"50\u00f72=25, 0\"; New = \"65\u00f74=16, 1\" },\n    @{ Row = 1; Col = 4; Old = \"73\u00f72=36, 1\"; New = \"79\u00f75=15, 4\" },\n    @{ Row = 1; Col = 5; Old = \"92\u00f72=46, 0\"; New = \"59\u00f72=29, 1\" },\n    @{ Row = 5; Col = 1; Old = \"67\u00f72=33, 1\"; New = \"41\u00f77=5, 6\" },\n    @{ Row = 5; Col = 2; Old = \"33\u00f76=5, 3\"; New = \"44\u00f74=11, 0\" },\n    @{ Row = 5; Col = 3; Old = \"37\u00f72=18, 1\"; New = \"13\u00f72=6, 1\" },\n    @{ Row = 5; Col = 4; Old = \"57\u00f76=9, 3\"; New = \"54\u00f78=6, 6\" },\n    @{ Row = 5; Col = 5; Old = \"34\u00f72=17, 0\"; New = \"94\u00f77=13, 3\" },\n    @{ Row = 9; Col = 1; Old = \"93\u00f76=15, 3\"; New = \"77\u00f73=25, 2\" },\n    @{ Row = 9; Col = 2; Old = \"35\u00f76=5, 5\"; New = \"39\u00f78=4, 7\" },\n    @{ Row = 9; Col = 3; Old = \"63\u00f74=15, 3\"; New = \"18\u00f73=6, 0\" },\n    @{ Row = 9; Col = 4; Old = \"77\u00f74=19, 1\"; New = \"41\u00f74=10, 1\" },\n    @{ Row = 9; Col = 5; Old = \"81\u00f77=11, 4\"; New = \"99\u00f78=12, 3\" },\n    @{ Row = 13; Col = 1; Old = \"80\u00f73=26, 2\"; New = \"90\u00f78=11, 2\" },\n    @{ Row = 13; Col = 2; Old = \"66\u00f72=33, 0\"; New = \"30\u00f76=5, 0\" },\n    @{ Row = 13; Col = 3; Old = \"92\u00f77=13, 1\"; New = \"19\u00f79=2, 1\" },\n    @{ Row = 13; Col = 4; Old = \"30\u00f74=7, 2\"; New = \"71\u00f77=10, 1\" },\n    @{ Row = 13; Col = 5; Old = \"90\u00f78=11, 2\"; New = \"61\u00f79=6, 7\" },\n    @{ Row = 17; Col = 1; Old = \"64\u00f76=10, 4\"; New = \"43\u00f75=8, 3\" },\n    @{ Row = 17; Col = 2; Old = \"43\u00f78=5, 3\"; New = \"66\u00f79=7, 3\" },\n    @{ Row = 17; Col = 3; Old = \"67\u00f78=8, 3\"; New = \"79\u00f73=26, 1\" },\n    @{ Row = 17; Col = 4; Old = \"66\u00f76=11, 0\"; New = \"46\u00f74=11, 2\" },\n    @{ Row = 17; Col = 5; Old = \"57\u00f73=19, 0\"; New = \"90\u00f74=22, 2\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $table.Cell($r.Row, $r.Col)\n    $current = $cell.Range.Text\n    $current = $current.Replace([char]13, \"\").Replace([char]7, \"\")\n    if ($current -ne $r.Old) {\n        throw \"Unexpected existing text at row $($r.Row) col $($r.Col): '$current' (expected '$($r.Old)')\"\n    }\n    $cell.Range.Text = $r.New\n}\n"}
